$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("D2").Value = "42.879.76"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.295.85"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'299.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'96.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.62%  "
$ws.Range("D7").Value = "'0.507"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").Value = "'33.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.59%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  -4.35%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "'16.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.78%  "
$ws.Range("D15").Value = "'6.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "2.654.34"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").Value = "2.304.60"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'0.805"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "42.790.61"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "'11.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").Value = "'6.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'67.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'235.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").Value = "'24.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("D29").Value = "'167.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'33.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").Value = "'9.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "'4.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("D35").Value = "'4.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'16.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "'0.0690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.101"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "'1.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("D44").Value = "1.993.97"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").Value = "'0.0280"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  -6.14%  "
$ws.Range("D48").Value = "'2.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").Value = "2.522.02"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "'52.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.71%  "
$ws.Range("D51").Value = "'4.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.17%  "
